# Update the "想去人数" (F column) values on the "展览" and "全部类型" sheets.
# Both sheets contain identical data, so the same row/value map applies to each.

$wb = $excel.ActiveWorkbook

$updates = @{
    2  = 365
    3  = 356
    4  = 1866
    7  = 191
    9  = 304
    11 = 4448
    12 = 22
    13 = 337
    14 = 1232
    15 = 521
    17 = 806
    19 = 433
    21 = 213
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Range("F$row").Value = $updates[$row]
    }
}
